$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 39717666
$ws.Range("C1").Value = 39617666
$ws.Range("D1").Value = "Было Katta Doimiy 40, Стало: Katta Doimiy 100"
$ws.Range("E1").Value = "2024-10-24 23:18:41"
